$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.473.27'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.15'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9981'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9984'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.96'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3766'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.50'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.45%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3581'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08193'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.232'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9981'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.34'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.523'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.336'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001227'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.635.72'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.33'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06956'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.734'
$ws.Range("D20").ClearFormats()

$ws.Range("E21").Value = '  -1.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9984'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.53'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.462.52'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("E25").Value = '  +2.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.125'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.18'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.93%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.49'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.183'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.05'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.818.11'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.32%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.742'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.80%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.093'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +7.51%  '

$ws.Range("B34").Value = 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.50'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.65%  '

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.033'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -10.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02764'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2493'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08769'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.011'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06987'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.61'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.79%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7003'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.56%  '

$ws.Range("E43").Value = '  -1.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.55'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6460'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.325'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9985'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.22%  '

$ws.Range("E48").Value = '  -0.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07936'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.57'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.181'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.84%  '
